# daily auto push: 2026-01-28 02:31 UTC
# Insert a new timestamped reading (2026/01/28, 水, 9, 201) right after the
# existing 2026/01/28 rows (row 737), which pushes every subsequent row
# down by one and grows the sheet from A1:D779 to A1:D780.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 738; everything at/below row 738
# (including the former row 738..779) shifts down to 739..780.
$ws.Rows.Item(738).Insert()

# Column A holds a date formatted as plain text (e.g. "2026/01/28") in this
# sheet, not a real Excel date. Writing that string straight into a
# General-formatted cell would make Excel auto-convert it into a date
# serial number, so mark the cell as Text first, write the value, then
# strip the temporary formatting back off so the cell matches its
# neighbours (no explicit style).
$ws.Range("A738").NumberFormat = "@"
$ws.Range("A738").Value = "2026/01/28"
$ws.Range("A738").ClearFormats()

$ws.Range("B738").Value = "水"
$ws.Range("C738").Value = 9
$ws.Range("D738").Value = 201
